$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Caballito Sur, Caballito'
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = '$ 600.000'
$ws.Cells.Item(2, 3).Value = '36 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(2, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-departamento-dos-ambientes-muy-luminoso-en-58010507.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=1'

$ws.Cells.Item(3, 1).Value = 'Balvanera, Capital Federal'
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = '$ 700.000'
$ws.Cells.Item(3, 3).Value = '41 m² tot.3 amb.2 dorm.1 baño'
$ws.Cells.Item(3, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-alquiler-3-ambientes-dueno-caba-balvanera-58080275.html?n_src=Listado&n_pills=Encargado&n_pg=1&n_pos=2'

$ws.Cells.Item(4, 1).Value = 'Almagro Sur, Almagro'
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = '$ 800.000'
$ws.Cells.Item(4, 3).Value = '42 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(4, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-2-amb-42-m-sup2--excelente-estado.-muy-58093867.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=3'

$ws.Cells.Item(5, 1).Value = 'Monserrat, Capital Federal'
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = '$ 1.400.000'
$ws.Cells.Item(5, 3).Value = '41 m² tot.3 amb.2 dorm.1 baño'
$ws.Cells.Item(5, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-de-3-ambientes-amoblado-listo-para-58096171.html?n_src=Listado&n_pills=Encargado&n_pg=1&n_pos=4'

$ws.Cells.Item(6, 1).Value = 'Tribunales, Capital Federal'
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = '$ 680.000'
$ws.Cells.Item(6, 3).Value = '37 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(6, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-depto-equipado-y-funcional-en-inmejorable-ubicacion!-52711772.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=5'

$ws.Cells.Item(7, 1).Value = 'Chacarita, Capital Federal'
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = '$ 500.000'
$ws.Cells.Item(7, 3).Value = '36 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(7, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-alquila-dpto-2-av-corrientes-y-dorrego-lateral-57241867.html?n_src=Listado&n_pills=Pileta&n_pg=1&n_pos=6'

$ws.Cells.Item(8, 1).Value = 'Palermo Soho, Palermo'
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = '$ 600.000'
$ws.Cells.Item(8, 3).Value = '65 m² tot.1 amb.1 dorm.1 baño1 coch.'
$ws.Cells.Item(8, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquilo-departamento-totalmente-amueblado-58096033.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=7'

$ws.Cells.Item(9, 1).Value = 'Caballito Sur, Caballito'
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = '$ 1.200.000'
$ws.Cells.Item(9, 3).Value = '83 m² tot.4 amb.3 dorm.1 baño'
$ws.Cells.Item(9, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-excelente-dpto-4-amb-c-dep-al-frente-58096013.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=8'

$ws.Cells.Item(10, 1).Value = 'Recoleta, Capital Federal'
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = '$ 450.000'
$ws.Cells.Item(10, 3).Value = '29 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(10, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-alquila-dpto-2-ambientes-lateral-super-luminoso-57250660.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=9'

$ws.Cells.Item(11, 1).Value = 'Almagro Norte, Almagro'
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = '$ 650.000'
$ws.Cells.Item(11, 3).Value = '38 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(11, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-alquila-dpto-2-amb-vista-abierta-super-luminoso-57652927.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=10'

$ws.Cells.Item(12, 1).Value = 'Villa del Parque, Capital Federal'
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = '$ 750.000'
$ws.Cells.Item(12, 3).Value = '41 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(12, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-con-terraza-y-sol-57938037.html?n_src=Listado&n_pills=Terraza&n_pg=1&n_pos=11'

$ws.Cells.Item(13, 1).Value = 'Belgrano R, Belgrano'
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = '$ 700.000'
$ws.Cells.Item(13, 3).Value = '60 m² tot.2 amb.1 dorm.1 baño1 coch.'
$ws.Cells.Item(13, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-en-belgrano-r-dueno-directo-58095968.html?n_src=Listado&n_pills=Pileta&n_pg=1&n_pos=12'

$ws.Cells.Item(14, 1).Value = 'Flores Norte, Flores'
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = '$ 620.000'
$ws.Cells.Item(14, 3).Value = '75 m² tot.3 amb.2 dorm.1 baño1 coch.'
$ws.Cells.Item(14, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095956.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=13'

$ws.Cells.Item(15, 1).Value = 'Flores Norte, Flores'
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = '$ 620.000'
$ws.Cells.Item(15, 3).Value = '75 m² tot.3 amb.2 dorm.1 baño1 coch.'
$ws.Cells.Item(15, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095955.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=14'

$ws.Cells.Item(16, 1).Value = 'Floresta Sur, Floresta'
$ws.Cells.Item(16, 2).NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = '$ 750'
$ws.Cells.Item(16, 3).Value = '89 m² tot.3 amb.2 dorm.2 baños1 coch.'
$ws.Cells.Item(16, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquilo-depto-amueblado-58095949.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=15'

$ws.Cells.Item(17, 1).Value = 'Villa Crespo, Capital Federal'
$ws.Cells.Item(17, 2).NumberFormat = "@"
$ws.Cells.Item(17, 2).Value = 'USD 800'
$ws.Cells.Item(17, 3).Value = '72 m² tot.3 amb.2 dorm.1 baño'
$ws.Cells.Item(17, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-dpto-3-ambientes-dueno-directo-amoblado-58095947.html?n_src=Listado&n_pills=Terraza&n_pg=1&n_pos=16'

$ws.Cells.Item(18, 1).Value = 'Flores Norte, Flores'
$ws.Cells.Item(18, 2).NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = '$ 620.000'
$ws.Cells.Item(18, 3).Value = '75 m² tot.3 amb.2 dorm.1 baño1 coch.'
$ws.Cells.Item(18, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095948.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=17'

$ws.Cells.Item(19, 1).Value = 'Flores Norte, Flores'
$ws.Cells.Item(19, 2).NumberFormat = "@"
$ws.Cells.Item(19, 2).Value = '$ 620.000'
$ws.Cells.Item(19, 3).Value = '75 m² tot.3 amb.2 dorm.1 baño1 coch.'
$ws.Cells.Item(19, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095946.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=18'

$ws.Cells.Item(20, 1).Value = 'Flores Norte, Flores'
$ws.Cells.Item(20, 2).NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = '$ 630.000'
$ws.Cells.Item(20, 3).Value = '70 m² tot.2 amb.1 dorm.1 baño1 coch.'
$ws.Cells.Item(20, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-2-ambiente-58095945.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=19'

$ws.Cells.Item(21, 1).Value = 'Flores Norte, Flores'
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = '$ 630.000'
$ws.Cells.Item(21, 3).Value = '75 m² tot.3 amb.2 dorm.1 baño1 coch.'
$ws.Cells.Item(21, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-2-ambiente-58095921.html?n_src=Listado&n_pills=Laundry&n_pg=1&n_pos=20'

$ws.Cells.Item(22, 1).Value = 'Flores Norte, Flores'
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = '$ 620.000'
$ws.Cells.Item(22, 3).Value = '75 m² tot.3 amb.2 dorm.1 baño1 coch.'
$ws.Cells.Item(22, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095925.html?n_src=Listado&n_pills=Laundry&n_pg=1&n_pos=21'

$ws.Cells.Item(23, 1).Value = 'Villa Urquiza, Capital Federal'
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = '$ 850.000'
$ws.Cells.Item(23, 3).Value = '48 m² tot.1 amb.1 dorm.1 baño'
$ws.Cells.Item(23, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-2-ambientes-en-alquiler-al-frente-villa-urquiza-dueno-58095610.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=22'

$ws.Cells.Item(24, 1).Value = 'San Cristobal, Capital Federal'
$ws.Cells.Item(24, 2).NumberFormat = "@"
$ws.Cells.Item(24, 2).Value = '$ 499.000'
$ws.Cells.Item(24, 3).Value = '37 m² tot.1 amb.1 baño'
$ws.Cells.Item(24, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-monoambiente.-bajas-expensas.-sin-comision.-54758980.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=23'

$ws.Cells.Item(25, 1).Value = 'Villa General Mitre, Capital Federal'
$ws.Cells.Item(25, 2).NumberFormat = "@"
$ws.Cells.Item(25, 2).Value = '$ 700.000'
$ws.Cells.Item(25, 3).Value = '39 m² tot.2 amb.1 baño'
$ws.Cells.Item(25, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-dos-ambientes-villa-general-mitre-58095085.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=24'

$ws.Cells.Item(26, 1).Value = 'Parque Rivadavia, Caballito'
$ws.Cells.Item(26, 2).NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = '$ 600.000'
$ws.Cells.Item(26, 3).Value = '42 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(26, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dos-ambientes-en-caballito.-58094336.html?n_src=Listado&n_pills=Encargado&n_pg=1&n_pos=25'

$ws.Cells.Item(27, 1).Value = 'San Cristobal, Capital Federal'
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = '$ 650.000'
$ws.Cells.Item(27, 3).Value = '35 m² tot.'
$ws.Cells.Item(27, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquilo-departamento-dos-ambientes-san-cristobal-58094313.html?n_src=Listado&n_pg=1&n_pos=26'

$ws.Cells.Item(28, 1).Value = 'Recoleta, Capital Federal'
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = 'USD 700'
$ws.Cells.Item(28, 3).Value = '200 m² tot.4 amb.3 dorm.2 baños'
$ws.Cells.Item(28, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-directo-todo-incluido-resuelvo-hoy.-video-58094288.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=27'

$ws.Cells.Item(29, 1).Value = 'San Cristobal, Capital Federal'
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = '$ 980.000'
$ws.Cells.Item(29, 3).Value = '74 m² tot.3 amb.2 dorm.1 baño'
$ws.Cells.Item(29, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-luminoso-dueno-directo-58094293.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=28'

$ws.Cells.Item(30, 1).Value = 'Belgrano, Capital Federal'
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = '$ 950.000'
$ws.Cells.Item(30, 3).Value = '52 m² tot.2 amb.1 dorm.1 baño1 coch.'
$ws.Cells.Item(30, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-duena-alquila-2-ambientes-con-cochera.-amplio-y-52962262.html?n_src=Listado&n_pills=SUM&n_pg=1&n_pos=29'

$ws.Cells.Item(31, 1).Value = 'Lomas de Núñez, Núñez'
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = 'USD 750'
$ws.Cells.Item(31, 3).Value = '53 m² tot.2 amb.1 dorm.1 baño'
$ws.Cells.Item(31, 4).Value = 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-departamento-de-2-ambientes-en-nunez-caba-58094208.html?n_src=Listado&n_pills=Encargado&n_pg=1&n_pos=30'
